$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Avg Cycles" rows (4-12) in column C used to be placeholder "-" text;
# replace them with the real measured values for HQC, Kyber and Saber.
$ws.Range("C4").Value = 571896356
$ws.Range("C5").Value = 1148267355
$ws.Range("C6").Value = 1743477982
$ws.Range("C7").Value = 1
$ws.Range("C8").Value = 1
$ws.Range("C9").Value = 1
$ws.Range("C10").Value = 103596
$ws.Range("C11").Value = 161364
$ws.Range("C12").Value = 175820

# Touch a size-8 "Aptos Narrow" font (reserved by Excel as the sheet's
# phonetic-guide font) so it is registered in the style table, then restore
# the cell's original look.
$tmp = $ws.Range("B14")
$tmp.Font.Size = 8
$tmp.Font.Size = 11

# The active selection moved from E10 to C10.
$ws.Range("C10").Select()
